$wb = $excel.ActiveWorkbook

# --- Sheet: Significant Components ---
$ws1 = $wb.Worksheets.Item("Significant Components")
$ws1.Range("C2").Value = "['QSERV' 'QESL' 'QNOHLTH' 'QEDLESHI' 'QEXTRCT' 'QHISPC' 'PPUNIT' 'QFHH'`n 'PERCAP']"
$ws1.Range("C4").Value = "['QAGEDEP' 'QSSBEN' 'MEDAGE']"

# --- Sheet: Loading Factors ---
$ws2 = $wb.Worksheets.Item("Loading Factors")
$ws2.Range("A2").Value = "QSERV"
$ws2.Range("B2").Value = 0.5788572269334414
$ws2.Range("C2").Value = 0.3591286421625365
$ws2.Range("D2").Value = -0.2252310192538391
$ws2.Range("E2").Value = -0.03396987931677808
$ws2.Range("F2").Value = 0.2788530180455505
$ws2.Range("A3").Value = "QESL"
$ws2.Range("B3").Value = 0.7948904984323639
$ws2.Range("C3").Value = 0.1613825334468211
$ws2.Range("D3").Value = -0.03005137998308928
$ws2.Range("E3").Value = -0.2422038143621637
$ws2.Range("F3").Value = 0.2090795683955667
$ws2.Range("A4").Value = "QNOHLTH"
$ws2.Range("B4").Value = 0.6835885831892797
$ws2.Range("C4").Value = 0.4274092513584713
$ws2.Range("D4").Value = -0.1180986227688646
$ws2.Range("E4").Value = -0.1193419859649904
$ws2.Range("F4").Value = 0.2783959678079707
$ws2.Range("A5").Value = "QEDLESHI"
$ws2.Range("B5").Value = 0.8742754545694222
$ws2.Range("C5").Value = 0.2170987387450964
$ws2.Range("D5").Value = -0.01552035325745195
$ws2.Range("E5").Value = -0.1115375102429765
$ws2.Range("F5").Value = 0.1928903334309309
$ws2.Range("A6").Value = "QEXTRCT"
$ws2.Range("B6").Value = 0.7638644823218199
$ws2.Range("C6").Value = 0.1521521475853751
$ws2.Range("D6").Value = 0.01442734521604606
$ws2.Range("E6").Value = -0.2421830635724205
$ws2.Range("F6").Value = 0.0980451600571053
$ws2.Range("A7").Value = "QHISPC"
$ws2.Range("B7").Value = 0.8302558099045889
$ws2.Range("C7").Value = 0.3391482215877205
$ws2.Range("D7").Value = -0.1355137705952958
$ws2.Range("E7").Value = -0.1303520320439234
$ws2.Range("F7").Value = 0.1010584356086445
$ws2.Range("A8").Value = "PPUNIT"
$ws2.Range("B8").Value = 0.7423686833126053
$ws2.Range("C8").Value = -0.02716542760737626
$ws2.Range("D8").Value = -0.1540540875389355
$ws2.Range("E8").Value = 0.05737345372566562
$ws2.Range("F8").Value = -0.4605803611641273
$ws2.Range("A9").Value = "QFHH"
$ws2.Range("B9").Value = 0.5649722670226905
$ws2.Range("C9").Value = 0.2979097563303481
$ws2.Range("D9").Value = -0.09661690089936902
$ws2.Range("E9").Value = 0.2627433105522369
$ws2.Range("F9").Value = -0.02662914703716459
$ws2.Range("A10").Value = "PERCAP"
$ws2.Range("B10").Value = 0.491374773462963
$ws2.Range("C10").Value = 0.7088158683189665
$ws2.Range("D10").Value = -0.276333280934828
$ws2.Range("E10").Value = 0.05584477791462314
$ws2.Range("F10").Value = 0.1848051844184266
$ws2.Range("A11").Value = "QRICH"
$ws2.Range("B11").Value = 0.213013272636517
$ws2.Range("C11").Value = 0.8688064563089002
$ws2.Range("D11").Value = -0.1838949632965971
$ws2.Range("E11").Value = -0.01368674942245885
$ws2.Range("F11").Value = 0.2842078608904167
$ws2.Range("A12").Value = "MDHSEVAL"
$ws2.Range("B12").Value = 0.3700960064485249
$ws2.Range("C12").Value = 0.7551164634472819
$ws2.Range("D12").Value = -0.0187335771015702
$ws2.Range("E12").Value = -0.01138785929196793
$ws2.Range("F12").Value = -0.01732185127475163
$ws2.Range("A13").Value = "QAGEDEP"
$ws2.Range("B13").Value = -0.04257251893742132
$ws2.Range("C13").Value = -0.1105675146675588
$ws2.Range("D13").Value = 0.6577462371544496
$ws2.Range("E13").Value = 0.6416531589319452
$ws2.Range("F13").Value = -0.1145417532589742
$ws2.Range("A14").Value = "QSSBEN"
$ws2.Range("B14").Value = 0.01958181058519956
$ws2.Range("C14").Value = -0.05546784723718445
$ws2.Range("D14").Value = 0.7724236636677263
$ws2.Range("E14").Value = 0.1368589223097694
$ws2.Range("F14").Value = -0.1456422712177872
$ws2.Range("A15").Value = "MEDAGE"
$ws2.Range("B15").Value = -0.310032971999737
$ws2.Range("C15").Value = -0.2422549013545611
$ws2.Range("D15").Value = 0.7910147337033275
$ws2.Range("E15").Value = -0.01282271536411952
$ws2.Range("F15").Value = -0.2761421321301501
$ws2.Range("A16").Value = "QFEMLBR"
$ws2.Range("B16").Value = -0.2392027557165198
$ws2.Range("C16").Value = 0.08115484564614343
$ws2.Range("D16").Value = -0.02961507201535967
$ws2.Range("E16").Value = 0.7870156598199922
$ws2.Range("F16").Value = 0.002700677784901289
$ws2.Range("A17").Value = "QFEMALE"
$ws2.Range("B17").Value = -0.04558114523530072
$ws2.Range("C17").Value = -0.04877177825685924
$ws2.Range("D17").Value = 0.1719743026450634
$ws2.Range("E17").Value = 0.8732591131431345
$ws2.Range("F17").Value = -0.02403664950012062
$ws2.Range("A18").Value = "QRENTER"
$ws2.Range("B18").Value = 0.007369445308590607
$ws2.Range("C18").Value = 0.2404828217006364
$ws2.Range("D18").Value = -0.4212014999083357
$ws2.Range("E18").Value = -0.09993109697197729
$ws2.Range("F18").Value = 0.7593165704367276
$ws2.Range("A19").Value = "QPOVTY"
$ws2.Range("B19").Value = 0.3702864779202081
$ws2.Range("C19").Value = 0.1444159373175141
$ws2.Range("D19").Value = -0.3834144468126657
$ws2.Range("E19").Value = 0.08162395691949784
$ws2.Range("F19").Value = 0.4744661448571343
$ws2.Range("A20").Value = "QNOAUTO"
$ws2.Range("B20").Value = 0.1603481989507422
$ws2.Range("C20").Value = 0.06253916239658538
$ws2.Range("D20").Value = -0.1037263512179677
$ws2.Range("E20").Value = -0.01653642180616662
$ws2.Range("F20").Value = 0.6401537827355265

# --- Sheet: All Refactor Variances ---
$ws3 = $wb.Worksheets.Item("All Refactor Variances")
$ws3.Range("B2").Value = 5.357711907063933
$ws3.Range("C2").Value = 2.502681462413677
$ws3.Range("D2").Value = 2.317329284213393
$ws3.Range("E2").Value = 2.200799336869337
$ws3.Range("F2").Value = 2.006259110707016
$ws3.Range("G2").Value = 1.600413605322627
$ws3.Range("H2").Value = 1.000035936442935
$ws3.Range("I2").Value = 4.773468753771055
$ws3.Range("J2").Value = 3.404064001175622
$ws3.Range("K2").Value = 2.235666416925869
$ws3.Range("L2").Value = 2.054347870287409
$ws3.Range("M2").Value = 2.034144212558734
$ws3.Range("N2").Value = 5.085980812112077
$ws3.Range("O2").Value = 2.605486864154571
$ws3.Range("P2").Value = 2.248589716232898
$ws3.Range("Q2").Value = 2.067416303071898
$ws3.Range("R2").Value = 1.906728077536314
$ws3.Range("B3").Value = 0.1984337743357012
$ws3.Range("C3").Value = 0.09269190601532136
$ws3.Range("D3").Value = 0.08582701052642197
$ws3.Range("E3").Value = 0.08151108655071618
$ws3.Range("F3").Value = 0.07430589298914875
$ws3.Range("G3").Value = 0.0592745779749121
$ws3.Range("H3").Value = 0.03703836801640502
$ws3.Range("I3").Value = 0.2273080358938598
$ws3.Range("J3").Value = 0.1620982857702677
$ws3.Range("K3").Value = 0.1064603055678985
$ws3.Range("L3").Value = 0.09782608906130517
$ws3.Range("M3").Value = 0.09686401012184448
$ws3.Range("N3").Value = 0.2676832006374777
$ws3.Range("O3").Value = 0.1371308875870827
$ws3.Range("P3").Value = 0.1183468271701525
$ws3.Range("Q3").Value = 0.1088113843722052
$ws3.Range("R3").Value = 0.1003541093440165
$ws3.Range("B4").Value = 0.1984337743357012
$ws3.Range("C4").Value = 0.2911256803510226
$ws3.Range("D4").Value = 0.3769526908774445
$ws3.Range("E4").Value = 0.4584637774281607
$ws3.Range("F4").Value = 0.5327696704173095
$ws3.Range("G4").Value = 0.5920442483922216
$ws3.Range("H4").Value = 0.6290826164086266
$ws3.Range("I4").Value = 0.2273080358938598
$ws3.Range("J4").Value = 0.3894063216641275
$ws3.Range("K4").Value = 0.495866627232026
$ws3.Range("L4").Value = 0.5936927162933312
$ws3.Range("M4").Value = 0.6905567264151757
$ws3.Range("N4").Value = 0.2676832006374777
$ws3.Range("O4").Value = 0.4048140882245604
$ws3.Range("P4").Value = 0.523160915394713
$ws3.Range("Q4").Value = 0.6319722997669182
$ws3.Range("R4").Value = 0.7323264091109347
$ws3.Range("B5").Value = 0.3154335681194641
$ws3.Range("C5").Value = 0.1473445674663381
$ws3.Range("D5").Value = 0.1364320174930286
$ws3.Range("E5").Value = 0.1295713542619494
$ws3.Range("F5").Value = 0.1181178609152389
$ws3.Range("G5").Value = 0.09422383710633285
$ws3.Range("H5").Value = 0.05887679463764803
$ws3.Range("I5").Value = 0.3291663482504374
$ws3.Range("J5").Value = 0.23473565540626
$ws3.Range("K5").Value = 0.1541659091796212
$ws3.Range("L5").Value = 0.1416626401847403
$ws3.Range("M5").Value = 0.1402694469789409
$ws3.Range("N5").Value = 0.3655244400682652
$ws3.Range("O5").Value = 0.1872537790267096
$ws3.Range("P5").Value = 0.16160393193225
$ws3.Range("Q5").Value = 0.1485831768709602
$ws3.Range("R5").Value = 0.137034672101815

# --- Sheet: Final Variances ---
$ws4 = $wb.Worksheets.Item("Final Variances")
$ws4.Range("B2").Value = 5.085980812112077
$ws4.Range("C2").Value = 2.605486864154571
$ws4.Range("D2").Value = 2.248589716232898
$ws4.Range("E2").Value = 2.067416303071898
$ws4.Range("F2").Value = 1.906728077536314
$ws4.Range("B3").Value = 0.2676832006374777
$ws4.Range("C3").Value = 0.1371308875870827
$ws4.Range("D3").Value = 0.1183468271701525
$ws4.Range("E3").Value = 0.1088113843722052
$ws4.Range("F3").Value = 0.1003541093440165
$ws4.Range("B4").Value = 0.2676832006374777
$ws4.Range("C4").Value = 0.4048140882245604
$ws4.Range("D4").Value = 0.523160915394713
$ws4.Range("E4").Value = 0.6319722997669182
$ws4.Range("F4").Value = 0.7323264091109347
$ws4.Range("B5").Value = 0.3655244400682652
$ws4.Range("C5").Value = 0.1872537790267096
$ws4.Range("D5").Value = 0.16160393193225
$ws4.Range("E5").Value = 0.1485831768709602
$ws4.Range("F5").Value = 0.137034672101815

# --- Sheet: Included and Excluded ---
$ws5 = $wb.Worksheets.Item("Included and Excluded")
$ws5.Range("B2").Value = "[['QSERV', 'QESL', 'QNOHLTH', 'QEDLESHI', 'QEXTRCT', 'QHISPC', 'PPUNIT', 'QFHH', 'PERCAP', 'QRICH', 'MDHSEVAL', 'QAGEDEP', 'QSSBEN', 'MEDAGE', 'QFEMLBR', 'QFEMALE', 'QRENTER', 'QPOVTY', 'QNOAUTO']]"
